$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold/border/center) from H1 onto the new headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data cells, row 2
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 6

# New data cells, row 3
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 7
